# Weekly update: add the newest week's price record for
# "Bruselas (repollito)" at Vega Central Mapocho de Santiago, and
# re-insert the previously-last record (which the source feed re-sends)
# as its own row at the end, per the author's "Fruta / hortaliza, semanal"
# refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new latest-week row at row 73 (pushes the existing
#     data rows 73-99 down to 74-100). Excel's row insert copies the
#     formatting (incl. the date number format on column D) from the
#     row above automatically.
$ws.Rows.Item(73).Insert()

$ws.Range("A73").Value = 9
$ws.Range("B73").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C73").Value = "Metropolitana"
$ws.Range("D73").Value = 45120
$ws.Range("E73").Value = 13
$ws.Range("F73").Value = 100112035
$ws.Range("G73").Value = "Bruselas (repollito)"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 34
$ws.Range("K73").Value = 18000
$ws.Range("L73").Value = 19000
$ws.Range("M73").Value = 18500
$ws.Range("N73").Value = "$/malla 15 kilos"
$ws.Range("O73").Value = "Provincia de Quillota"
$ws.Range("P73").Value = 1233
$ws.Range("Q73").Value = 15
$ws.Range("R73").Value = "Hortaliza"

# --- Insert a second new row. After the first insert, the old row 98
#     now lives at row 99, so inserting here pushes old-98 -> 100 and
#     old-99 -> 101, and row 99 becomes the slot for the new record.
$ws.Rows.Item(99).Insert()

$ws.Range("A99").Value = 9
$ws.Range("B99").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C99").Value = "Metropolitana"
$ws.Range("D99").Value = 45121
$ws.Range("E99").Value = 13
$ws.Range("F99").Value = 100112035
$ws.Range("G99").Value = "Bruselas (repollito)"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 52
$ws.Range("K99").Value = 18000
$ws.Range("L99").Value = 19000
$ws.Range("M99").Value = 18500
$ws.Range("N99").Value = "$/malla 15 kilos"
$ws.Range("O99").Value = "Provincia de Quillota"
$ws.Range("P99").Value = 1233
$ws.Range("Q99").Value = 15
$ws.Range("R99").Value = "Hortaliza"
